## Switch license from BY-NC to BY-SA (slide 2 of the Petascale lesson 11.3 deck)
##
## This script reproduces, via the PowerPoint COM/object model, the edit that:
##   1. Nudges the title placeholder's horizontal offset by 2 EMU
##      (566059 -> 566057) -- an incidental side effect of PowerPoint
##      re-flowing the text box when its text was edited.
##   2. Replaces "BY-NC" with "BY-SA" in the license statement, and
##      "by-nc" with "by-sa" in the Creative Commons hyperlink URL/text,
##      exactly like a user selecting those substrings and retyping them
##      (which is why PowerPoint ends up splitting the original single
##      runs into multiple runs around the edited text).
##   3. Removes the (empty / unused) <p:timing> animation-timing block
##      that had been left over on the slide.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$sh = $s.Shapes.Item(1)

# --- 1. Tiny re-flow nudge on the title placeholder's position ---------
# 566057 EMU / 12700 EMU-per-point = 44.57141732283465 pt
$sh.Left = 44.57141732283465

# --- 2. Swap the license from CC BY-NC 4.0 to CC BY-SA 4.0 -------------
$tf = $sh.TextFrame
$tr = $tf.TextRange
$fullText = $tr.Text

# a) "CC BY-NC 4.0. To view a copy of this license, visit "
#    -> only "BY-NC " changes to "BY-SA "
$ncIdx = $fullText.IndexOf("BY-NC ")
$ncRange = $tr.Characters($ncIdx + 1, 6)
$ncRange.Text = "BY-SA "

# b) hyperlink text/URL: https://creativecommons.org/licenses/by-nc/4.0
#    -> https://creativecommons.org/licenses/by-sa/4.0
#    (the user retyped everything after "https://")
$fullText = $tr.Text
$urlIdx = $fullText.IndexOf("creativecommons.org/licenses/by-nc/4.0")
$urlRange = $tr.Characters($urlIdx + 1, 39)
$urlRange.Text = "creativecommons.org/licenses/by-sa/4.0"

# --- 3. Drop the leftover empty animation timing on this slide ---------
# Adding then immediately deleting a throwaway effect clears the
# now-empty <p:timing> node from the slide's XML entirely.
$tmpEffect = $s.TimeLine.MainSequence.AddEffect($sh, 1)
$s.TimeLine.MainSequence.Item(1).Delete()
